$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.226.63"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").Value = "3.532.78"
$ws.Range("E3").Value = "  -0.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.56"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.34"
$ws.Range("E6").Value = "  -1.59%  "

# Row 7
$ws.Range("D7").Value = "3.532.48"
$ws.Range("E7").Value = "  -0.17%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"

# Row 10
$ws.Range("E10").Value = "  +0.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  +3.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  +0.23%  "

# Row 13
$ws.Range("D13").Value = "4.150.46"
$ws.Range("E13").Value = "  +0.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.54"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15
$ws.Range("E15").Value = "  +0.16%  "

# Row 16
$ws.Range("D16").Value = "3.543.64"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("D18").Value = "65.324.87"
$ws.Range("E18").Value = "  +0.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.80"
$ws.Range("E19").Value = "  -3.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.52"
$ws.Range("E20").Value = "  +2.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.75"
$ws.Range("E21").Value = "  -2.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.43"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("E23").Value = "  +1.30%  "

# Row 24
$ws.Range("D24").Value = "3.679.68"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.94"
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  +2.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.91"
$ws.Range("E28").Value = "  +0.83%  "

# Row 29
$ws.Range("E29").Value = "  +15.37%  "

# Row 30
$ws.Range("E30").Value = "  +2.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.30"
$ws.Range("E31").Value = "  +1.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.46"
$ws.Range("E32").Value = "  +2.33%  "

# Row 33
$ws.Range("D33").Value = "3.541.41"
$ws.Range("E33").Value = "  -0.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.20"
$ws.Range("E34").Value = "  +1.34%  "

# Row 35
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.146"
$ws.Range("E36").Value = "  +1.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.32"
$ws.Range("E37").Value = "  +5.93%  "

# Row 38
$ws.Range("E38").Value = "  -0.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.52"
$ws.Range("E39").Value = "  +0.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.90"
$ws.Range("E40").Value = "  -0.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0829"
$ws.Range("E41").Value = "  +3.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.826"
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.06"
$ws.Range("E43").Value = "  -2.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.25"
$ws.Range("E44").Value = "  +4.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.81"
$ws.Range("E45").Value = "  +0.83%  "

# Row 46
$ws.Range("E46").Value = "  +0.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.47"
$ws.Range("E47").Value = "  +0.50%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.67"
$ws.Range("E48").Value = "  -0.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.92"
$ws.Range("E49").Value = "  +1.07%  "

# Row 50
$ws.Range("D50").Value = "2.384.91"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51 - VeChain -> SuiNetwork
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.897"
$ws.Range("E51").Value = "  +5.84%  "
